$d = $word.ActiveDocument

$replacements = @(
    @{old = "88÷3=29, 1"; new = "43÷4=10, 3"},
    @{old = "52÷7=7, 3"; new = "81÷4=20, 1"},
    @{old = "87÷4=21, 3"; new = "80÷4=20, 0"},
    @{old = "18÷8=2, 2"; new = "80÷4=20, 0"},
    @{old = "68÷7=9, 5"; new = "40÷7=5, 5"},
    @{old = "14÷9=1, 5"; new = "81÷4=20, 1"},
    @{old = "43÷6=7, 1"; new = "20÷3=6, 2"},
    @{old = "16÷8=2, 0"; new = "95÷2=47, 1"},
    @{old = "67÷4=16, 3"; new = "59÷7=8, 3"},
    @{old = "60÷6=10, 0"; new = "52÷8=6, 4"},
    @{old = "98÷2=49, 0"; new = "74÷7=10, 4"},
    @{old = "41÷2=20, 1"; new = "11÷6=1, 5"},
    @{old = "49÷8=6, 1"; new = "47÷2=23, 1"},
    @{old = "77÷2=38, 1"; new = "26÷9=2, 8"},
    @{old = "63÷2=31, 1"; new = "32÷7=4, 4"},
    @{old = "46÷3=15, 1"; new = "43÷2=21, 1"},
    @{old = "44÷3=14, 2"; new = "38÷7=5, 3"},
    @{old = "99÷6=16, 3"; new = "68÷8=8, 4"},
    @{old = "94÷7=13, 3"; new = "48÷7=6, 6"},
    @{old = "52÷3=17, 1"; new = "49÷9=5, 4"},
    @{old = "63÷9=7, 0"; new = "41÷3=13, 2"},
    @{old = "46÷2=23, 0"; new = "81÷5=16, 1"},
    @{old = "60÷9=6, 6"; new = "43÷8=5, 3"},
    @{old = "91÷5=18, 1"; new = "11÷3=3, 2"},
    @{old = "63÷7=9, 0"; new = "75÷5=15, 0"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
